# Fill in Oil (O), Protein (P), and Oil+Protein (Q) values for HIF 6
# sheet rows that previously had these cells blank. Q = O + P (lsmeans
# calculated for oil/protein, used to clean up/complete the test-weight
# yield export).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HIF 6")

$data = @(
    @{ Row = 2; O = 24.14; P = 40.92; Q = 65.06 },
    @{ Row = 3; O = 24.03; P = 41.09; Q = 65.12 },
    @{ Row = 5; O = 22.23; P = 44.55; Q = 66.78 },
    @{ Row = 6; O = 22.11; P = 43.93; Q = 66.04 },
    @{ Row = 8; O = 22.68; P = 42.34; Q = 65.02 },
    @{ Row = 9; O = 23.06; P = 41.24; Q = 64.3 },
    @{ Row = 11; O = 24.8; P = 38.48; Q = 63.28 },
    @{ Row = 12; O = 24.04; P = 40.07; Q = 64.11 },
    @{ Row = 15; O = 22.08; P = 44.78; Q = 66.86 },
    @{ Row = 16; O = 21.68; P = 45.2; Q = 66.88 },
    @{ Row = 18; O = 21.34; P = 46.7; Q = 68.04 },
    @{ Row = 19; O = 20.93; P = 47.6; Q = 68.53 },
    @{ Row = 21; O = 24.81; P = 39.31; Q = 64.12 },
    @{ Row = 22; O = 24.54; P = 40.77; Q = 65.31 },
    @{ Row = 23; O = 24.54; P = 40.77; Q = 65.31 },
    @{ Row = 25; O = 24.41; P = 41.23; Q = 65.64 },
    @{ Row = 26; O = 24.41; P = 41.23; Q = 65.64 },
    @{ Row = 27; O = 23.92; P = 42.04; Q = 65.96 },
    @{ Row = 28; O = 23.92; P = 42.04; Q = 65.96 },
    @{ Row = 30; O = 24.05; P = 39.64; Q = 63.69 },
    @{ Row = 31; O = 24.11; P = 41.3; Q = 65.41 },
    @{ Row = 32; O = 24.11; P = 41.3; Q = 65.41 },
    @{ Row = 35; O = 24.46; P = 39.98; Q = 64.44 },
    @{ Row = 36; O = 23.49; P = 41.07; Q = 64.56 },
    @{ Row = 38; O = 23.69; P = 40.83; Q = 64.52 },
    @{ Row = 39; O = 24.6; P = 39.29; Q = 63.89 },
    @{ Row = 41; O = 25.21; P = 38.07; Q = 63.28 },
    @{ Row = 42; O = 24.04; P = 40.34; Q = 64.38 },
    @{ Row = 44; O = 23.9; P = 39.87; Q = 63.77 },
    @{ Row = 45; O = 24.45; P = 39.83; Q = 64.28 },
    @{ Row = 47; O = 25.06; P = 39.42; Q = 64.48 },
    @{ Row = 48; O = 24.88; P = 38.94; Q = 63.82 },
    @{ Row = 50; O = 21.7; P = 44.36; Q = 66.06 },
    @{ Row = 51; O = 22.54; P = 43.18; Q = 65.72 },
    @{ Row = 53; O = 23.72; P = 42.58; Q = 66.3 },
    @{ Row = 54; O = 23; P = 44.29; Q = 67.29 },
    @{ Row = 56; O = 20.96; P = 47.63; Q = 68.59 },
    @{ Row = 57; O = 20.44; P = 47.75; Q = 68.19 },
    @{ Row = 59; O = 21.02; P = 47.45; Q = 68.47 },
    @{ Row = 60; O = 20.73; P = 49.01; Q = 69.74 },
    @{ Row = 62; O = 20.48; P = 47.02; Q = 67.5 },
    @{ Row = 63; O = 20.03; P = 48.83; Q = 68.86 },
    @{ Row = 65; O = 20.57; P = 48.05; Q = 68.62 },
    @{ Row = 66; O = 20.53; P = 48.44; Q = 68.97 },
    @{ Row = 68; O = 24.64; P = 39.74; Q = 64.38 },
    @{ Row = 69; O = 23.79; P = 41.2; Q = 64.99 },
    @{ Row = 71; O = 24.15; P = 40.92; Q = 65.07 },
    @{ Row = 72; O = 23.48; P = 41.22; Q = 64.7 },
    @{ Row = 74; O = 21.86; P = 44.71; Q = 66.57 },
    @{ Row = 75; O = 21.51; P = 45.57; Q = 67.08 },
    @{ Row = 77; O = 22.29; P = 44.91; Q = 67.2 },
    @{ Row = 78; O = 22.07; P = 46.28; Q = 68.35 },
    @{ Row = 80; O = 25.29; P = 40.69; Q = 65.98 },
    @{ Row = 81; O = 25.07; P = 40.97; Q = 66.04 },
    @{ Row = 83; O = 23.38; P = 43.65; Q = 67.03 },
    @{ Row = 84; O = 22.41; P = 44.65; Q = 67.06 },
    @{ Row = 86; O = 21.38; P = 46.3; Q = 67.68 },
    @{ Row = 87; O = 21.62; P = 45.87; Q = 67.49 },
    @{ Row = 89; O = 22.83; P = 39.41; Q = 62.24 },
    @{ Row = 90; O = 22.46; P = 40.32; Q = 62.78 },
    @{ Row = 92; O = 21.5; P = 44.23; Q = 65.73 },
    @{ Row = 93; O = 21.88; P = 43.49; Q = 65.37 },
    @{ Row = 95; O = 23.45; P = 39.16; Q = 62.61 },
    @{ Row = 96; O = 22.93; P = 38.89; Q = 61.82 }
)

foreach ($item in $data) {
    $row = $item.Row
    $ws.Cells.Item($row, 15).Value = $item.O   # Column O - Oil
    $ws.Cells.Item($row, 16).Value = $item.P   # Column P - Protein
    $ws.Cells.Item($row, 17).Value = $item.Q   # Column Q - Oil + Protein
}
